# "updated PvsI with model fitting"
# Adds an "effective_volume" column (G) to the MasterSheet, computed from
# the Buoyant_Weight(g) column (F) via the standard buoyant-weight ->
# displaced-volume conversion: (154.4 - (F / 1.025)) / 1000

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterSheet")

# Header for the new column
$ws.Range("G1").Value = "effective_volume"

# Formulas for rows 2-23 (Excel will detect/emit these as a shared formula
# group, same as the existing column E formulas)
$ws.Range("G2:G23").Formula = "=(154.4 - (F2 / 1.025)) / 1000"

# Leave the new selection where the author's session ended up
$ws.Range("L21").Select()
